$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill out the date series across row 3 and row 4 (pandas-style
#     navigation structure: continue the existing day-by-day date run
#     that row 2 already has, out to column U) ---

$cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U")

# Row 3 continues from E3 (42043) => F3=42044 ... U3=42059
$row3Start = 42044
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "3").Value = $row3Start + $i
}

# Row 4 continues from E4 (42044) => F4=42045 ... U4=42060
$row4Start = 42045
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "4").Value = $row4Start + $i
}

# Match the date number formatting already used in the row (copy format
# from the last populated cell of each row across the newly filled cells)
$ws.Range("E3").Copy()
$ws.Range("F3:U3").PasteSpecial(-4122)
$ws.Range("E4").Copy()
$ws.Range("F4:U4").PasteSpecial(-4122)

# --- Update the sheet's navigation/view state: move the selection to D5 ---
$ws.Range("D5").Select()
